# Ran through presentation and timed it.
# Swap the order of slide 2 ("BlaTeX") and slide 3 ("Background") so that
# the "Background" slide now comes right after the title slide, followed
# by the "BlaTeX" slide.

$p = $ppt.ActivePresentation

# Slide 3 ("Background") moves up to position 2; slide 2 ("BlaTeX") is
# pushed down to position 3.
$p.Slides.Item(3).MoveTo(2)
